# The "ID" sheet maps acquisition channel ids -> friendly channel names.
# This edit removes the fluorescence channel rows (flu1 / flu2) from the
# map, per "removed some channel maps for better coverage / want some
# channels to not be in the channel map".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ID")
$ws.Activate()

# Row 4 (485530 -> flu1): drop the mapped value but keep the row/style
# so A4 remains formatted (numFmt) yet empty.
$ws.Range("A4:B4").ClearContents()

# Row 5 (485,530[2] -> flu2): remove the row entirely.
$ws.Rows.Item(5).Delete()

# Reflect the saved selection state from the edited workbook.
$ws.Range("B7").Select()
